# Add an "upload essay quiz" sheet: rename the existing sheet to
# MULTIPLE_CHOICE, duplicate it as ESSAY, trim the essay sheet down to
# NAME/VALUE/INFO columns only, and fill in the essay question/answer data.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "MULTIPLE_CHOICE"

# Tag the original sheet with its quiz type.
$ws1.Range("C2").Value2 = "MULTIPLE_CHOICE"

# Duplicate the whole sheet (keeps styles/merges/column widths identical)
# and place the copy right after the source sheet.
$ws1.Copy([System.Type]::Missing, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "ESSAY"

# The essay sheet only needs the first three columns (NAME / VALUE / INFO).
$ws2.Range("D1:I14").EntireColumn.Delete()

# Essay-specific column widths.
$ws2.Columns.Item(2).ColumnWidth = 26.09
$ws2.Columns.Item(3).ColumnWidth = 55.26

# Quiz-type + answer-column labels.
$ws2.Range("C2").Value2 = "ESSAY"
$ws2.Range("C5").Value2 = "ANSWER ESSAY"

# Essay questions/answers (replace the MULTIPLE_CHOICE formulas with literal
# text, and fill in the per-row answer column).
for ($i = 7; $i -le 14; $i++) {
    $n = $i - 6
    $ws2.Cells.Item($i, 2).Value2 = "question essay $n"
    $ws2.Cells.Item($i, 3).Value2 = "Answer essay $n"
}

# Restore the original sheet's selection/view and select the essay sheet's
# question range, matching the finished workbook layout.
$ws1.Range("A1:C14").Select()
$ws2.Range("B7:B14").Select()
$ws2.Activate()
